$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 178659
$ws.Range("C4").Value = 168615
$ws.Range("C7").Value = 5.62
$ws.Range("C8").Value = 65.05
